$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph "Play Alien Antix Free - Exciting Gameplay with
#    Cluster Wins".
# ---------------------------------------------------------------------------
$firstHeading = $d.Paragraphs.Item(1)
$firstHeading.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Borrow the run layout (leading empty run + one formatted run) of the bold
# "Play Alien Antix Free..." paragraph near the end of the document so the
# freshly created paragraph gets the same <w:r/> + <w:r> shape used
# throughout the rest of the document.
$boldTemplatePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$metaPara.Range.FormattedText = $boldTemplatePara.Range.FormattedText

# Re-fetch the paragraph (positions were just rewritten) and swap its text
# for "Meta description" while keeping the Bold run formatting.
$metaPara = $d.Paragraphs.Item(2)
$boldRun = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$boldRun.Text = "Meta description"

# Append the remaining (non-bold) text of the meta description right before
# the paragraph mark.
$metaPara = $d.Paragraphs.Item(2)
$pilcrow = $metaPara.Range.End - 1
$tail = $d.Range($pilcrow, $pilcrow)
$tail.InsertBefore(": Read our Alien Antix review for exciting gameplay with bonuses, high payouts, and cluster wins. Play Alien Antix for free now.")

$metaPara = $d.Paragraphs.Item(2)
$tailRange = $d.Range($pilcrow, $metaPara.Range.End - 1)
$tailRange.Bold = 0

# ---------------------------------------------------------------------------
# 2) Remove the old bold "Play Alien Antix Free..." paragraph that used to
#    sit near the end of the document (right before the italic meta
#    description paragraph).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$oldBoldPara = $d.Paragraphs.Item($count - 1)
$oldBoldPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the trailing italic paragraph (formerly the "Read
#    our Alien Antix review..." meta description) with the new AI image
#    generation prompt, keeping the italic formatting intact.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$italicRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$italicRange.Text = 'Create an eye-catching feature image for the online slot game "Alien Antix". The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. Make sure to incorporate elements of space and aliens in the image to match the theme of the game. The image should be vibrant and colorful, with the Maya warrior as the central focus, surrounded by aliens and other space objects. Make the image stand out to attract potential players to the game.'
